$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.410.55"
$ws.Range("E2").Value = "  +8.82%  "
$ws.Range("D3").Value = "1.603.07"
$ws.Range("E3").Value = "  +8.11%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("E4").Value = "  -0.56%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.9920"
$ws.Range("E5").Value = "  +2.17%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "301.51"
$ws.Range("E6").Value = "  +7.58%  "
$ws.Range("E7").Value = "  +0.79%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3406"
$ws.Range("E8").Value = "  +9.81%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "42.57"
$ws.Range("E9").Value = "  +5.89%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.141"
$ws.Range("E10").Value = "  +6.93%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07078"
$ws.Range("E11").Value = "  +5.76%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.002"
$ws.Range("E12").Value = "  -0.41%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.77"
$ws.Range("E13").Value = "  +8.66%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.948"
$ws.Range("E14").Value = "  +7.18%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.648"
$ws.Range("E15").Value = "  +6.68%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001091"
$ws.Range("D17").Value = "1.602.49"
$ws.Range("E17").Value = "  +8.05%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9922"
$ws.Range("E18").Value = "  +2.22%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06834"
$ws.Range("E19").Value = "  +14.56%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "77.95"
$ws.Range("E20").Value = "  +11.39%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.048"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "16.15"
$ws.Range("E22").Value = "  +10.68%  "
$ws.Range("E23").Value = "  +6.94%  "
$ws.Range("D24").Value = "22.450.03"
$ws.Range("E24").Value = "  +8.61%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.400"
$ws.Range("E25").Value = "  +5.79%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.541"
$ws.Range("E26").Value = "  +19.49%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "150.92"
$ws.Range("E27").Value = "  +6.06%  "
$ws.Range("E28").Value = "  +12.95%  "
$ws.Range("D29").Value = "1.782.36"
$ws.Range("E29").Value = "  +8.51%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "121.08"
$ws.Range("E30").Value = "  +5.29%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.171"
$ws.Range("E31").Value = "  +5.34%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.114"
$ws.Range("E32").Value = "  +20.88%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.9542"
$ws.Range("E33").Value = "  +15.43%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08277"
$ws.Range("E34").Value = "  +3.03%  "
$ws.Range("E35").Value = "  +6.07%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.300"
$ws.Range("E36").Value = "  +10.91%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "11.95"
$ws.Range("E37").Value = "  +14.00%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.275"
$ws.Range("E38").Value = "  +4.03%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.604"
$ws.Range("E39").Value = "  +11.80%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06120"
$ws.Range("E40").Value = "  +5.20%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.02223"
$ws.Range("E41").Value = "  +8.33%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.2029"
$ws.Range("E42").Value = "  +7.54%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9915"
$ws.Range("E43").Value = "  +2.18%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5926"
$ws.Range("E44").Value = "  +11.15%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.812"
$ws.Range("E45").Value = "  +7.60%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "13.21"
$ws.Range("E46").Value = "  +7.18%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5703"
$ws.Range("E47").Value = "  +9.13%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "127.85"
$ws.Range("E48").Value = "  +7.43%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.974"
$ws.Range("E49").Value = "  +7.90%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06817"
$ws.Range("E50").Value = "  +4.69%  "
$ws.Range("E51").Value = "  +8.74%  "
